$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for columns D (Fecha, as Excel date serials), J (Volumen),
# K (Precio minimo), L (Precio maximo), M (Precio promedio ponderado) and
# P (Precio $/Kg) for rows 3-12. These values represent the rows' data
# being re-ordered (the underlying records got reshuffled / re-dated).

$rows = @{
    3  = @{ D = 44832; J = 60; K = 17000; L = 18000; M = 17500; P = 1346 }
    4  = @{ D = 44841; J = 30; K = 18000; L = 18000; M = 18000; P = 1385 }
    5  = @{ D = 44846; J = 30; K = 18000; L = 18000; M = 18000; P = 1385 }
    6  = @{ D = 44880; J = 30; K = 17000; L = 17000; M = 17000; P = 1308 }
    7  = @{ D = 44804; J = 40; K = 12000; L = 13000; M = 12500; P = 962 }
    8  = @{ D = 44859; J = 30; K = 13000; L = 13000; M = 13000; P = 1000 }
    9  = @{ D = 44797; J = 60; K = 12000; L = 13000; M = 12500; P = 962 }
    10 = @{ D = 44868; J = 30; K = 18000; L = 18000; M = 18000; P = 1385 }
    11 = @{ D = 44874; J = 30; K = 17000; L = 17000; M = 17000; P = 1308 }
    12 = @{ D = 44810; J = 40; K = 12000; L = 13000; M = 12500; P = 962 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value = [DateTime]::FromOADate($vals.D)
    $ws.Range("J$r").Value = $vals.J
    $ws.Range("K$r").Value = $vals.K
    $ws.Range("L$r").Value = $vals.L
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("P$r").Value = $vals.P
}
